$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sciences")

# Rows where category code C moves from 9 -> 10 (keeping the same fractional part in F)
$upRows = @(5, 59, 64, 65)
foreach ($r in $upRows) {
    $cCell = $ws.Cells.Item($r, 3)   # column C
    $fCell = $ws.Cells.Item($r, 6)   # column F
    $cCell.Value2 = $cCell.Value2 + 1
    $fCell.Value2 = $fCell.Value2 + 1
}

# Rows where category code C moves from 10 -> 9
$downRows = @(68, 69)
foreach ($r in $downRows) {
    $cCell = $ws.Cells.Item($r, 3)   # column C
    $fCell = $ws.Cells.Item($r, 6)   # column F
    $cCell.Value2 = $cCell.Value2 - 1
    $fCell.Value2 = $fCell.Value2 - 1
}

$wb.Application.Calculate()
